$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 09:52"

# --- Row 18: Rusia - refreshed case counts ---
$ws.Range("B18").Value = 21102
$ws.Range("C18").Value = 2774
$ws.Range("D18").Value = 1694
$ws.Range("E18").Value = 19238
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 170

# --- Row 20: Austria - refreshed case counts ---
$ws.Range("B20").Value = 14083
$ws.Range("C20").Value = 42
$ws.Range("E20").Value = 6372

# --- Row 33: Australia - refreshed case counts ---
$ws.Range("B33").Value = 6400
$ws.Range("C33").Value = 41
$ws.Range("D33").Value = 3598
$ws.Range("E33").Value = 2741
$ws.Range("F33").Value = 80

# --- Rows 65-67: Estonia overtakes Nueva Zelanda and Barein in ranking ---
# Row 65 becomes Estonia with refreshed data
$ws.Range("A65").Value = "Estonia"
$ws.Range("B65").Value = 1373
$ws.Range("C65").Value = 41
$ws.Range("D65").Value = 115
$ws.Range("E65").Value = 1227
$ws.Range("F65").Value = 11
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 31

# Row 66 becomes Nueva Zelanda (previous row 65 data, unchanged)
$ws.Range("A66").Value = "Nueva Zelanda"
$ws.Range("B66").Value = 1366
$ws.Range("C66").Value = 17
$ws.Range("D66").Value = 628
$ws.Range("E66").Value = 729
$ws.Range("F66").Value = 4
$ws.Range("G66").Value = 4
$ws.Range("H66").Value = 9

# Row 67 becomes Barein (previous row 66 data, unchanged)
$ws.Range("A67").Value = "Barein"
$ws.Range("B67").Value = 1361
$ws.Range("C67").Value = 0
$ws.Range("D67").Value = 591
$ws.Range("E67").Value = 763
$ws.Range("F67").Value = 3
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 7

# --- Rows 79-81: Eslovaquia overtakes Oman and Banglades in ranking ---
# Row 79 becomes Eslovaquia with refreshed data
$ws.Range("A79").Value = "Eslovaquia"
$ws.Range("B79").Value = 816
$ws.Range("C79").Value = 47
$ws.Range("D79").Value = 107
$ws.Range("E79").Value = 707
$ws.Range("F79").Value = 5
$ws.Range("H79").Value = 2

# Row 80 becomes Oman (previous row 79 data, unchanged)
$ws.Range("A80").Value = "Oman"
$ws.Range("B80").Value = 813
$ws.Range("C80").Value = 86
$ws.Range("D80").Value = 130
$ws.Range("E80").Value = 679
$ws.Range("F80").Value = 3
$ws.Range("H80").Value = 4

# Row 81 becomes Banglades (previous row 80 data, unchanged)
$ws.Range("A81").Value = "Banglades"
$ws.Range("B81").Value = 803
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 42
$ws.Range("E81").Value = 722
$ws.Range("F81").Value = 1
$ws.Range("H81").Value = 39

# --- Row 99: Honduras - refreshed case counts ---
$ws.Range("B99").Value = 407
$ws.Range("C99").Value = 10
$ws.Range("E99").Value = 374
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 26

# --- Row 114: Vietnam - refreshed case counts ---
$ws.Range("D114").Value = 166
$ws.Range("E114").Value = 99
